$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the values in columns D, J, K, L, M, P between row 3 and row 4
$cols = @("D", "J", "K", "L", "M", "P")

foreach ($col in $cols) {
    $cell3 = $ws.Range($col + "3")
    $cell4 = $ws.Range($col + "4")
    $temp = $cell3.Value2
    $cell3.Value2 = $cell4.Value2
    $cell4.Value2 = $temp
}
